$d = $word.ActiveDocument

# NOTE on technique: mutating text through Find.Execute's Replace
# argument (or through a Range obtained from Cell.Range / Paragraph.Range
# etc.) in this runtime is not reliably confined to the Range it is
# called on - it can locate/replace unrelated, identically-formatted text
# far away in the document. To stay surgical we therefore only ever
# resolve *positions* with Find (on Ranges freshly built with
# $d.Range(start, end), which *is* properly bounded) and perform the
# actual edits with plain Range.Text assignment, which only ever touches
# the characters inside the given Range.

function Find-InBounds($startPos, $endPos, $text) {
    $r = $d.Range($startPos, $endPos)
    $ok = $r.Find.Execute($text, $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
    if (-not $ok) {
        throw "Text not found: $text"
    }
    return $r
}

# =========================================================================
# Change 1: consolidate the "...этапа системы" / ". Плакат - фо" / "р" /
# "мат А1, лист 1." run split back into "...этапа системы" (unchanged,
# its own run) + ". Плакат - формат А1, лист 1." (one merged run).
# =========================================================================

$docEnd = $d.Content.End

# "контент-аналитического" only occurs once, right before the affected
# text, so it is a reliable, unique anchor.
$anchor = Find-InBounds 0 $docEnd "контент-аналитического"

# Bound the local search window to just after the anchor so Find cannot
# stray into any of the other, look-alike "Плакат - формат А1, лист 1."
# captions elsewhere in the document.
$windowEnd = [Math]::Min($anchor.End + 200, $docEnd)

$etapa = Find-InBounds $anchor.End $windowEnd " этапа системы"
$splitPlakatText = ". Плакат - фо" + "р" + "мат А1, лист 1."
$plakatSpan = Find-InBounds $etapa.End $windowEnd $splitPlakatText

# Overwrite the three split runs' combined span with the same text as one
# write - this merges them (and only them, because the write is confined
# to exactly their span) into a single run.
$mergeRange = $d.Range($plakatSpan.Start, $plakatSpan.End)
$mergeRange.Text = ". Плакат - формат А1, лист 1."

# =========================================================================
# Change 2: "7. Содержание задания по экологической безопасности" ->
# "7. Содержание задания по охране труда", keeping the leading phrase in
# its own run and putting the new tail in a second run.
# =========================================================================

$docEnd2 = $d.Content.End
$titleSpan = Find-InBounds 0 $docEnd2 "7. Содержание задания по экологической безопасности"
$prefixSpan = Find-InBounds $titleSpan.Start $titleSpan.End "7. Содержание задания по "

$newSuffixText = "охране труда"
$suffixRange = $d.Range($prefixSpan.End, $titleSpan.End)
$suffixRange.Text = $newSuffixText

# Force the new tail text onto its own run (identical character
# formatting would otherwise make the engine re-coalesce it with the
# preceding, unchanged prefix run) by round-tripping a character-format
# property over just the new span - this does not change its visible
# formatting since the property is toggled back immediately.
$newSuffixLen = $newSuffixText.Length
$newSuffixEnd = $prefixSpan.End + $newSuffixLen
$newSuffixRange = $d.Range($prefixSpan.End, $newSuffixEnd)
$newSuffixRange.Font.Bold = $true
$newSuffixRange.Font.Bold = $false
